$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 46
$ws.Range("H46").Value = 4016.6667
$ws.Range("J46").Value = 4275
$ws.Range("L46").Value = 12825
$ws.Range("N46").Value = -13063

# Row 60
$ws.Range("H60").Value = 4016.6667
$ws.Range("J60").Value = 4275
$ws.Range("L60").Value = 12825
$ws.Range("N60").Value = -13793

# Row 64
$ws.Range("H64").Value = 3809.25
$ws.Range("J64").Value = 4666.6665
$ws.Range("L64").Value = 4666.6665
$ws.Range("N64").Value = -5162.6665

# Row 67
$ws.Range("H67").Value = 3809.25
$ws.Range("J67").Value = 4666.6665
$ws.Range("L67").Value = 4666.6665
$ws.Range("N67").Value = -6382.6665

# Row 76
$ws.Range("H76").Value = 3871.4285
$ws.Range("I76").Value = 3730.7693
$ws.Range("J76").Value = 4277.778
$ws.Range("K76").Value = 3730.7693
$ws.Range("L76").Value = 4277.778
$ws.Range("M76").Value = -3415.7693
$ws.Range("N76").Value = -4907.778

# Row 79
$ws.Range("H79").Value = 3871.4285
$ws.Range("I79").Value = 3730.7693
$ws.Range("J79").Value = 4277.778
$ws.Range("K79").Value = 3730.7693
$ws.Range("L79").Value = 4277.778
$ws.Range("M79").Value = -2638.7693
$ws.Range("N79").Value = -6461.778

# Row 80
$ws.Range("H80").Value = 988.5714
$ws.Range("I80").Value = 415.18182
$ws.Range("J80").Value = 1619.3
$ws.Range("K80").Value = 1245.54546
$ws.Range("L80").Value = 4857.9
$ws.Range("M80").Value = -247.54546
$ws.Range("N80").Value = -6853.9

# Row 83
$ws.Range("H83").Value = 988.5714
$ws.Range("I83").Value = 415.18182
$ws.Range("J83").Value = 1619.3
$ws.Range("K83").Value = 3736.63638
$ws.Range("L83").Value = 14573.7
$ws.Range("M83").Value = 1255.36362
$ws.Range("N83").Value = -24557.7

# Row 137
$ws.Range("H137").Value = 3634.9412
$ws.Range("I137").Value = 3786.2666
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 11358.7998
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -8808.799800000001
$ws.Range("N137").Value = -12600

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1379.3611
$ws.Range("I2").Value = 1375.9143
$ws.Range("K2").Value = 1375.9143
$ws.Range("M2").Value = -1262.9143

# Row 32
$ws.Range("H32").Value = 33874.695
$ws.Range("I32").Value = 33874.695
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 33874.695
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -33587.695
$ws.Range("N32").ClearContents()

# Row 45
$ws.Range("H45").Value = 1604.2174
$ws.Range("J45").Value = 1266.3334
$ws.Range("L45").Value = 1266.3334
$ws.Range("N45").Value = -2020.3334

# Row 63
$ws.Range("H63").Value = 3074.4443
$ws.Range("I63").Value = 2525.7144
$ws.Range("J63").Value = 4995
$ws.Range("K63").Value = 2525.7144
$ws.Range("L63").Value = 4995
$ws.Range("M63").Value = -1839.7144
$ws.Range("N63").Value = -6367

# Row 66
$ws.Range("H66").Value = 3074.4443
$ws.Range("I66").Value = 2525.7144
$ws.Range("J66").Value = 4995
$ws.Range("K66").Value = 12628.572
$ws.Range("L66").Value = 24975
$ws.Range("M66").Value = -9196.572
$ws.Range("N66").Value = -31839

# Row 74
$ws.Range("H74").Value = 2437.5667
$ws.Range("I74").Value = 2082.2856
$ws.Range("J74").Value = 3266.5557
$ws.Range("K74").Value = 2082.2856
$ws.Range("L74").Value = 3266.5557
$ws.Range("M74").Value = -1208.2856
$ws.Range("N74").Value = -5014.5557

# Row 77
$ws.Range("H77").Value = 2437.5667
$ws.Range("I77").Value = 2082.2856
$ws.Range("J77").Value = 3266.5557
$ws.Range("K77").Value = 10411.428
$ws.Range("L77").Value = 16332.7785
$ws.Range("M77").Value = -6043.428
$ws.Range("N77").Value = -25068.7785

# Row 116
$ws.Range("H116").Value = 1379.3611
$ws.Range("I116").Value = 1375.9143
$ws.Range("K116").Value = 1375.9143
$ws.Range("M116").Value = 918.0857000000001

# Row 132
$ws.Range("H132").Value = 6295.577
$ws.Range("I132").Value = 6851.8945
$ws.Range("J132").Value = 4785.5713
$ws.Range("K132").Value = 20555.6835
$ws.Range("L132").Value = 14356.7139
$ws.Range("M132").Value = -18025.6835
$ws.Range("N132").Value = -19416.7139

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1379.3611
$ws.Range("I3").Value = 1375.9143
$ws.Range("K3").Value = 1375.9143
$ws.Range("M3").Value = -1261.9143

# Row 20
$ws.Range("H20").Value = 1757.875
$ws.Range("I20").Value = 2193.8572
$ws.Range("J20").Value = 1418.7778
$ws.Range("K20").Value = 2193.8572
$ws.Range("L20").Value = 1418.7778
$ws.Range("M20").Value = -1946.8572
$ws.Range("N20").Value = -1912.7778

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 10485.083
$ws.Range("I16").Value = 17970.166
$ws.Range("K16").Value = 17970.166
$ws.Range("M16").Value = -17683.166

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# Row 113
$ws.Range("H113").Value = 10485.083
$ws.Range("I113").Value = 17970.166
$ws.Range("K113").Value = 17970.166
$ws.Range("M113").Value = -15800.166

# Row 122
$ws.Range("H122").Value = 3928.8096
$ws.Range("I122").Value = 3650.2778
$ws.Range("K122").Value = 10950.8334
$ws.Range("M122").Value = -8500.8334

# Row 132
$ws.Range("H132").Value = 2563.75
$ws.Range("I132").Value = 2125.5
$ws.Range("J132").Value = 3352.6
$ws.Range("K132").Value = 6376.5
$ws.Range("L132").Value = 10057.8
$ws.Range("M132").Value = -3846.5
$ws.Range("N132").Value = -15117.8

# Row 134
$ws.Range("H134").Value = 2451.5
$ws.Range("I134").Value = 2239.4358
$ws.Range("J134").Value = 4105.6
$ws.Range("K134").Value = 6718.307400000001
$ws.Range("L134").Value = 12316.8
$ws.Range("M134").Value = -4183.307400000001
$ws.Range("N134").Value = -17386.8

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 34.86207
$ws.Range("I2").Value = 22.75
$ws.Range("J2").Value = 36.8
$ws.Range("K2").Value = 136.5
$ws.Range("L2").Value = 220.8
$ws.Range("M2").Value = -23.5
$ws.Range("N2").Value = -446.8

# Row 22
$ws.Range("H22").Value = 2133.3333
$ws.Range("J22").Value = 2360
$ws.Range("L22").Value = 7080
$ws.Range("N22").Value = -7418

# Row 27
$ws.Range("H27").Value = 2133.3333
$ws.Range("J27").Value = 2360
$ws.Range("L27").Value = 7080
$ws.Range("N27").Value = -7284

# Row 43
$ws.Range("H43").Value = 4674.5
$ws.Range("J43").Value = 4674.5
$ws.Range("L43").Value = 14023.5
$ws.Range("N43").Value = -14251.5

# Row 61
$ws.Range("H61").Value = 371.6
$ws.Range("I61").Value = 76.5
$ws.Range("J61").Value = 568.3333
$ws.Range("K61").Value = 229.5
$ws.Range("L61").Value = 1704.9999
$ws.Range("M61").Value = -14.5
$ws.Range("N61").Value = -2134.9999

# Row 81
$ws.Range("H81").Value = 4149.3335
$ws.Range("I81").Value = 2526.5
$ws.Range("K81").Value = 7579.5
$ws.Range("M81").Value = -6456.5

# Row 84
$ws.Range("H84").Value = 4149.3335
$ws.Range("I84").Value = 2526.5
$ws.Range("K84").Value = 22738.5
$ws.Range("M84").Value = -17122.5

# Row 131
$ws.Range("H131").Value = 15970.381
$ws.Range("I131").Value = 1110.0714
$ws.Range("J131").Value = 20216.184
$ws.Range("K131").Value = 3330.2142
$ws.Range("L131").Value = 60648.552
$ws.Range("M131").Value = 1709.7858
$ws.Range("N131").Value = -70728.552

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5509.154
$ws.Range("I70").Value = 5236.185
$ws.Range("J70").Value = 6123.3335
$ws.Range("K70").Value = 5236.185
$ws.Range("L70").Value = 6123.3335
$ws.Range("M70").Value = -4966.185
$ws.Range("N70").Value = -6663.3335

# Row 73
$ws.Range("H73").Value = 5509.154
$ws.Range("I73").Value = 5236.185
$ws.Range("J73").Value = 6123.3335
$ws.Range("K73").Value = 5236.185
$ws.Range("L73").Value = 6123.3335
$ws.Range("M73").Value = -4300.185
$ws.Range("N73").Value = -7995.3335

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2900
$ws.Range("I93").Value = 2825
$ws.Range("K93").Value = 2825
$ws.Range("M93").Value = -1577

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 3343.3333
$ws.Range("I107").Value = 679.3333
$ws.Range("J107").Value = 11335.333
$ws.Range("K107").Value = 2037.9999
$ws.Range("L107").Value = 34005.999
$ws.Range("M107").Value = -117.9999
$ws.Range("N107").Value = -37845.999

# Row 141
$ws.Range("H141").Value = 45000
$ws.Range("J141").Value = 45000
$ws.Range("L141").Value = 45000
$ws.Range("N141").Value = -55360
